# Generate Report for Archive
# - Flip the "Ready for handoff" status to "In Translation" everywhere it is
#   used (Overview!E2/F2 summary columns + the per-locale Status column on
#   the "zh-cn" and "de-de" sheets).
# - Narrow the (now shorter) status columns: the "17.22"-wide columns that
#   held "Ready for handoff" shrink to fit "In Translation".

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

# --- Update the status text everywhere it appears ---
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$zhcn.Range("C2").Value = $newStatus
$dede.Range("C2").Value = $newStatus

# --- Narrow the columns that used to hold the status text ---
$newWidth = 12.5

$overview.Columns.Item(5).ColumnWidth = $newWidth
$overview.Columns.Item(6).ColumnWidth = $newWidth
$zhcn.Columns.Item(3).ColumnWidth = $newWidth
$dede.Columns.Item(3).ColumnWidth = $newWidth
